$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the second row entirely (A2:C2) so only row 1 remains
$ws.Rows.Item(2).Delete()

# Update B1 with the new long text
$ws.Range("B1").Value = "انما الاعمال بالنيات وانما لكل امرئ ما نوى فمن كانت هجرته الى الله و رسوله فهجرته الى الله و رسوله و من كانت هجرته لدنيا يصيبها او امراه ينكحها فهجرته الى ما هاجر اليه"

# Set column B width
$ws.Columns.Item(2).ColumnWidth = 92.1666666666667

# Set row 1 height
$ws.Rows.Item(1).RowHeight = 162

# Update the selection to B1
$ws.Range("B1").Select()
